$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 7356535
$ws.Range("I41").Value = 11367314
$ws.Range("K41").Value = 11367314
$ws.Range("M41").Value = -11366874
$ws.Range("H53").Value = 5146.615
$ws.Range("I53").Value = 4655.5713
$ws.Range("J53").Value = 5719.5
$ws.Range("K53").Value = 4655.5713
$ws.Range("L53").Value = 5719.5
$ws.Range("M53").Value = -4018.5713
$ws.Range("N53").Value = -6993.5
$ws.Range("H70").Value = 36448656
$ws.Range("I70").Value = 14449684
$ws.Range("J70").Value = 69447110
$ws.Range("K70").Value = 43349052
$ws.Range("L70").Value = 208341330
$ws.Range("M70").Value = -43348782
$ws.Range("N70").Value = -208341870
$ws.Range("H73").Value = 36448656
$ws.Range("I73").Value = 14449684
$ws.Range("J73").Value = 69447110
$ws.Range("K73").Value = 43349052
$ws.Range("L73").Value = 208341330
$ws.Range("M73").Value = -43348116
$ws.Range("N73").Value = -208343202
$ws.Range("H74").Value = 107151530
$ws.Range("I74").Value = 214289000
$ws.Range("K74").Value = 214289000
$ws.Range("M74").Value = -214288064
$ws.Range("H76").Value = 25253750
$ws.Range("I76").Value = 25253750
$ws.Range("K76").Value = 25253750
$ws.Range("M76").Value = -25253435
$ws.Range("H77").Value = 107151530
$ws.Range("I77").Value = 214289000
$ws.Range("K77").Value = 1071445000
$ws.Range("M77").Value = -1071440320
$ws.Range("H79").Value = 25253750
$ws.Range("I79").Value = 25253750
$ws.Range("K79").Value = 25253750
$ws.Range("M79").Value = -25252658
$ws.Range("H111").Value = 17861278
$ws.Range("I111").Value = 41667990
$ws.Range("K111").Value = 125003970
$ws.Range("M111").Value = -125000903
$ws.Range("H132").Value = 1716.119
$ws.Range("I132").Value = 1693.5588
$ws.Range("J132").Value = 1812
$ws.Range("K132").Value = 5080.6764
$ws.Range("L132").Value = 5436
$ws.Range("M132").Value = -2550.6764
$ws.Range("N132").Value = -10496
$ws.Range("H137").Value = 3903.8262
$ws.Range("I137").Value = 4971.5713
$ws.Range("J137").Value = 3436.6875
$ws.Range("K137").Value = 14914.7139
$ws.Range("L137").Value = 10310.0625
$ws.Range("M137").Value = -12364.7139
$ws.Range("N137").Value = -15410.0625
$ws.Range("H138").Value = 1542033.4
$ws.Range("I138").Value = 677.4857
$ws.Range("J138").Value = 3340282
$ws.Range("K138").Value = 2032.4571
$ws.Range("L138").Value = 10020846
$ws.Range("M138").Value = 3107.5429
$ws.Range("N138").Value = -10031126
$ws.Range("H141").Value = 1463.4286
$ws.Range("I141").Value = 875.5
$ws.Range("K141").Value = 2626.5
$ws.Range("M141").Value = 2553.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1990667.5
$ws.Range("I32").Value = 2277222.5
$ws.Range("J32").Value = 20601.5
$ws.Range("K32").Value = 2277222.5
$ws.Range("L32").Value = 20601.5
$ws.Range("M32").Value = -2276935.5
$ws.Range("N32").Value = -21175.5
$ws.Range("H61").Value = 9965
$ws.Range("I61").Value = 2342.1428
$ws.Range("K61").Value = 2342.1428
$ws.Range("M61").Value = -2130.1428
$ws.Range("H64").Value = 30627
$ws.Range("J64").Value = 33999.5
$ws.Range("L64").Value = 33999.5
$ws.Range("N64").Value = -34495.5
$ws.Range("H67").Value = 30627
$ws.Range("J67").Value = 33999.5
$ws.Range("L67").Value = 33999.5
$ws.Range("N67").Value = -35715.5
$ws.Range("H74").Value = 30568.182
$ws.Range("I74").Value = 42593.332
$ws.Range("K74").Value = 42593.332
$ws.Range("M74").Value = -41719.332
$ws.Range("H77").Value = 30568.182
$ws.Range("I77").Value = 42593.332
$ws.Range("K77").Value = 212966.66
$ws.Range("M77").Value = -208598.66
$ws.Range("H96").Value = 46971
$ws.Range("J96").Value = 46971
$ws.Range("L96").Value = 46971
$ws.Range("N96").Value = -52463
$ws.Range("H98").Value = 54544.332
$ws.Range("J98").Value = 54544.332
$ws.Range("L98").Value = 54544.332
$ws.Range("N98").Value = -60534.332
$ws.Range("H132").Value = 4332.855
$ws.Range("I132").Value = 2118.9768
$ws.Range("J132").Value = 9343.210999999999
$ws.Range("K132").Value = 6356.930399999999
$ws.Range("L132").Value = 28029.633
$ws.Range("M132").Value = -3826.930399999999
$ws.Range("N132").Value = -33089.633
$ws.Range("H136").Value = 9965
$ws.Range("I136").Value = 2342.1428
$ws.Range("K136").Value = 7026.428400000001
$ws.Range("M136").Value = -4476.428400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 6412309
$ws.Range("I20").Value = 8773930
$ws.Range("J20").Value = 2196.2856
$ws.Range("K20").Value = 8773930
$ws.Range("L20").Value = 2196.2856
$ws.Range("M20").Value = -8773683
$ws.Range("N20").Value = -2690.2856
$ws.Range("H80").Value = 29412090
$ws.Range("J80").Value = 291.9091
$ws.Range("L80").Value = 291.9091
$ws.Range("N80").Value = -2287.9091
$ws.Range("H83").Value = 29412090
$ws.Range("J83").Value = 291.9091
$ws.Range("L83").Value = 1459.5455
$ws.Range("N83").Value = -11443.5455
$ws.Range("H94").Value = 1093.4359
$ws.Range("I94").Value = 871.10345
$ws.Range("J94").Value = 1738.2
$ws.Range("K94").Value = 871.10345
$ws.Range("L94").Value = 1738.2
$ws.Range("M94").Value = -420.10345
$ws.Range("N94").Value = -2640.2
$ws.Range("H105").Value = 4598.5557
$ws.Range("I105").Value = 3999.5
$ws.Range("K105").Value = 3999.5
$ws.Range("M105").Value = -2252.5
$ws.Range("H116").Value = 70000
$ws.Range("J116").Value = 70000
$ws.Range("L116").Value = 70000
$ws.Range("N116").Value = -79178
$ws.Range("H134").Value = 5166.396
$ws.Range("I134").Value = 1885.5883
$ws.Range("K134").Value = 5656.7649
$ws.Range("M134").Value = -3121.7649

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 34232.6
$ws.Range("J28").Value = 34232.6
$ws.Range("L28").Value = 34232.6
$ws.Range("N28").Value = -34722.6
$ws.Range("H58").Value = 6761181.5
$ws.Range("I58").Value = 10640221
$ws.Range("J58").Value = 8780.333000000001
$ws.Range("K58").Value = 10640221
$ws.Range("L58").Value = 8780.333000000001
$ws.Range("M58").Value = -10640018
$ws.Range("N58").Value = -9186.333000000001
$ws.Range("H99").Value = 7942.1055
$ws.Range("J99").Value = 6858.5835
$ws.Range("L99").Value = 6858.5835
$ws.Range("N99").Value = -9854.583500000001
$ws.Range("H105").Value = 4468142
$ws.Range("I105").Value = 6495021.5
$ws.Range("K105").Value = 6495021.5
$ws.Range("M105").Value = -6493274.5
$ws.Range("H107").Value = 2233.1924
$ws.Range("I107").Value = 2124.4546
$ws.Range("K107").Value = 2124.4546
$ws.Range("M107").Value = -204.4546
$ws.Range("H126").Value = 7942.1055
$ws.Range("J126").Value = 6858.5835
$ws.Range("L126").Value = 20575.7505
$ws.Range("N126").Value = -25515.7505
$ws.Range("H132").Value = 5198898
$ws.Range("I132").Value = 2135.9607
$ws.Range("K132").Value = 6407.882100000001
$ws.Range("M132").Value = -3877.882100000001
$ws.Range("H134").Value = 4942.4165
$ws.Range("I134").Value = 2685.3
$ws.Range("J134").Value = 7763.8125
$ws.Range("K134").Value = 8055.900000000001
$ws.Range("L134").Value = 23291.4375
$ws.Range("M134").Value = -5520.900000000001
$ws.Range("N134").Value = -28361.4375
$ws.Range("H136").Value = 6761181.5
$ws.Range("I136").Value = 10640221
$ws.Range("J136").Value = 8780.333000000001
$ws.Range("K136").Value = 31920663
$ws.Range("L136").Value = 26340.999
$ws.Range("M136").Value = -31918113
$ws.Range("N136").Value = -31440.999
$ws.Range("H141").Value = 126981.664
$ws.Range("J141").Value = 126981.664
$ws.Range("L141").Value = 126981.664
$ws.Range("N141").Value = -137341.664

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 145.14285
$ws.Range("I6").Value = 145.14285
$ws.Range("K6").Value = 435.42855
$ws.Range("M6").Value = -322.42855
$ws.Range("H23").Value = 319.73334
$ws.Range("I23").Value = 247.66667
$ws.Range("J23").Value = 367.77777
$ws.Range("K23").Value = 743.00001
$ws.Range("L23").Value = 1103.33331
$ws.Range("M23").Value = -508.00001
$ws.Range("N23").Value = -1573.33331
$ws.Range("H95").Value = 10000
$ws.Range("J95").Value = 10000
$ws.Range("L95").Value = 30000
$ws.Range("N95").Value = -34118
$ws.Range("H113").Value = 3357.318
$ws.Range("J113").Value = 4702.4165
$ws.Range("L113").Value = 14107.2495
$ws.Range("N113").Value = -18447.2495
$ws.Range("H132").Value = 8156.826
$ws.Range("I132").Value = 3850.25
$ws.Range("J132").Value = 18000.428
$ws.Range("K132").Value = 34652.25
$ws.Range("L132").Value = 162003.852
$ws.Range("M132").Value = -32122.25
$ws.Range("N132").Value = -167063.852
$ws.Range("H133").Value = 8399.4
$ws.Range("I133").Value = 7333
$ws.Range("J133").Value = 9999
$ws.Range("K133").Value = 21999
$ws.Range("L133").Value = 29997
$ws.Range("M133").Value = -16939
$ws.Range("N133").Value = -40117
$ws.Range("H134").Value = 5095.579
$ws.Range("I134").Value = 4523.3125
$ws.Range("K134").Value = 13569.9375
$ws.Range("M134").Value = -8499.9375
$ws.Range("H137").Value = 44509.875
$ws.Range("I137").Value = 1940.8572
$ws.Range("J137").Value = 62038.293
$ws.Range("K137").Value = 5822.571599999999
$ws.Range("L137").Value = 186114.879
$ws.Range("M137").Value = -722.5715999999993
$ws.Range("N137").Value = -196314.879
$ws.Range("H138").Value = 6920
$ws.Range("I138").Value = 5508.5
$ws.Range("K138").Value = 16525.5
$ws.Range("M138").Value = -11385.5
$ws.Range("H139").Value = 4542.025
$ws.Range("I139").Value = 2387.28
$ws.Range("K139").Value = 7161.84
$ws.Range("M139").Value = -2021.84
$ws.Range("H140").Value = 4074.8
$ws.Range("I140").Value = 1855.3334
$ws.Range("K140").Value = 5566.0002
$ws.Range("M140").Value = -386.0002000000004
$ws.Range("H141").Value = 5330.577
$ws.Range("I141").Value = 2314.0476
$ws.Range("K141").Value = 6942.1428
$ws.Range("M141").Value = -1762.1428

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 55565960
$ws.Range("I70").Value = 142865900
$ws.Range("K70").Value = 142865900
$ws.Range("M70").Value = -142865630
$ws.Range("H73").Value = 55565960
$ws.Range("I73").Value = 142865900
$ws.Range("K73").Value = 142865900
$ws.Range("M73").Value = -142864964
$ws.Range("H80").Value = 80363.84
$ws.Range("I80").Value = 3466.375
$ws.Range("J80").Value = 203399.8
$ws.Range("K80").Value = 3466.375
$ws.Range("L80").Value = 203399.8
$ws.Range("M80").Value = -2468.375
$ws.Range("N80").Value = -205395.8
$ws.Range("H83").Value = 80363.84
$ws.Range("I83").Value = 3466.375
$ws.Range("J83").Value = 203399.8
$ws.Range("K83").Value = 17331.875
$ws.Range("L83").Value = 1016999
$ws.Range("M83").Value = -12339.875
$ws.Range("N83").Value = -1026983
$ws.Range("H102").Value = 1490.2632
$ws.Range("I102").Value = 1263.7446
$ws.Range("K102").Value = 1263.7446
$ws.Range("M102").Value = 358.2554
$ws.Range("H113").Value = 5942.7393
$ws.Range("I113").Value = 2666.7727
$ws.Range("J113").Value = 8945.708000000001
$ws.Range("K113").Value = 2666.7727
$ws.Range("L113").Value = 8945.708000000001
$ws.Range("M113").Value = -496.7727
$ws.Range("N113").Value = -13285.708
$ws.Range("H122").Value = 4480183
$ws.Range("I122").Value = 10234753
$ws.Range("J122").Value = 4406.6665
$ws.Range("K122").Value = 30704259
$ws.Range("L122").Value = 13219.9995
$ws.Range("M122").Value = -30701809
$ws.Range("N122").Value = -18119.9995
$ws.Range("H124").Value = 75540.5
$ws.Range("J124").Value = 75540.5
$ws.Range("L124").Value = 75540.5
$ws.Range("N124").Value = -85360.5
$ws.Range("H126").Value = 6557.5713
$ws.Range("I126").Value = 6123.1113
$ws.Range("J126").Value = 7339.6
$ws.Range("K126").Value = 18369.3339
$ws.Range("L126").Value = 22018.8
$ws.Range("M126").Value = -15899.3339
$ws.Range("N126").Value = -26958.8
$ws.Range("H132").Value = 5742
$ws.Range("I132").Value = 2296.6667
$ws.Range("J132").Value = 11307.538
$ws.Range("K132").Value = 6890.000100000001
$ws.Range("L132").Value = 33922.614
$ws.Range("M132").Value = -4360.000100000001
$ws.Range("N132").Value = -38982.614

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 308.05264
$ws.Range("I16").Value = 308.05264
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 308.05264
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -138.05264
$ws.Range("N16").ClearContents()
$ws.Range("H43").Value = 7500
$ws.Range("I43").Value = 7000
$ws.Range("K43").Value = 7000
$ws.Range("M43").Value = -6807
$ws.Range("H55").Value = 257.88635
$ws.Range("I55").Value = 91.208336
$ws.Range("K55").Value = 91.208336
$ws.Range("M55").Value = 81.791664
$ws.Range("H95").Value = 200000
$ws.Range("J95").Value = 200000
$ws.Range("L95").Value = 200000
$ws.Range("N95").Value = -205492
$ws.Range("H101").Value = 41894.375
$ws.Range("J101").Value = 41894.375
$ws.Range("L101").Value = 41894.375
$ws.Range("N101").Value = -48384.375
$ws.Range("H132").Value = 7251691.5
$ws.Range("I132").Value = 13891243
$ws.Range("J132").Value = 8543.546
$ws.Range("K132").Value = 41673729
$ws.Range("L132").Value = 25630.638
$ws.Range("M132").Value = -41671199
$ws.Range("N132").Value = -30690.638
$ws.Range("H136").Value = 8480.243
$ws.Range("I136").Value = 2451.2666
$ws.Range("K136").Value = 7353.7998
$ws.Range("M136").Value = -4803.7998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 20883.857
$ws.Range("J18").Value = 20883.857
$ws.Range("L18").Value = 20883.857
$ws.Range("N18").Value = -21229.857
$ws.Range("H96").Value = 1538.8
$ws.Range("J96").Value = 1999.5
$ws.Range("L96").Value = 1999.5
$ws.Range("N96").Value = -4745.5
$ws.Range("H97").Value = 40000
$ws.Range("J97").Value = 40000
$ws.Range("L97").Value = 40000
$ws.Range("N97").Value = -41982
$ws.Range("H122").Value = 163324.2
$ws.Range("J122").Value = 3791.5454
$ws.Range("L122").Value = 11374.6362
$ws.Range("N122").Value = -16274.6362
$ws.Range("H126").Value = 3552.4
$ws.Range("I126").Value = 1929
$ws.Range("J126").Value = 3958.25
$ws.Range("K126").Value = 5787
$ws.Range("L126").Value = 11874.75
$ws.Range("M126").Value = -3317
$ws.Range("N126").Value = -16814.75
$ws.Range("H132").Value = 12509802
$ws.Range("I132").Value = 15154991
$ws.Range("J132").Value = 39628
$ws.Range("K132").Value = 45464973
$ws.Range("L132").Value = 118884
$ws.Range("M132").Value = -45462443
$ws.Range("N132").Value = -123944
$ws.Range("H136").Value = 27810434
$ws.Range("I136").Value = 58824596
$ws.Range("K136").Value = 176473788
$ws.Range("M136").Value = -176471238

